$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(6,20,20,20,1,2.5,1,0.5,0,15,3,0,0,0,0.01,0.01,-3.171132837305141,3.880973849795271,-8.498550820234328,0,"blink"),
    @(6,20,20,20,1,2.5,1,0.5,0,15,3,0,0,0,0.01,0.01,-4.778879412017504,-1.657886467100935,-14.30660717678537,1,"templerun"),
    @(6,20,20,20,1,2.5,1,0.5,0,15,3,0,0,0,0.01,0.01,-1.816601255302154,3.933759619071374,-11.43069734112945,2,"sudoku"),
    @(6,20,20,20,1,2.5,1,0.5,0,15,3,0,0,0,0.01,0.01,-2.224943738798284,8.485458883586197,-10.39795834166682,3,"flow"),
    @(6,20,20,20,1,2.5,1,0.5,0,15,3,0,0,0,0.01,0.01,-1.317635029509511,13.04536256885818,-9.736668748976236,4,"wordsearch")
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
